$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("startup")
$ws.Columns.Item(3).Insert()

$statQuery = "MATCH (s:study) WITH COLLECT(DISTINCT(s.clinical_study_designation)) AS all_studies MATCH (d:demographic) WITH COLLECT(DISTINCT(d.breed)) AS all_breeds, COLLECT(DISTINCT(d.sex)) AS all_sexes, all_studies MATCH (d:diagnosis) WITH COLLECT(DISTINCT(d.disease_term)) AS all_diseases, all_breeds, all_sexes, all_studies MATCH (p:program)<-[*]-(s:study)<-[*]-(c:case)<--(demo:demographic), (c)<--(diag:diagnosis) WHERE demo.breed IN['Akita']  OPTIONAL MATCH (f:file)-[*]->(c), (samp:sample)-[*]->(c) WITH DISTINCT c AS c, p, s, demo, diag, f, samp RETURN count(DISTINCT(f)) as number_of_files , count(DISTINCT(samp)) as number_of_sample , count(DISTINCT(c.case_id)) as number_of_cases , count(DISTINCT(s.clinical_study_designation)) as number_of_study"
$ws.Range("C2").Value = $statQuery
$ws.Range("C1").Value = "StatQuery"

$ws.Columns.Item(3).AutoFit()
$ws.Columns.Item(4).AutoFit()

Write-Host "C width:" $ws.Columns.Item(3).ColumnWidth
Write-Host "D width:" $ws.Columns.Item(4).ColumnWidth
